$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = "aa"
$ws.Range("J11").Value = "Agree/Accept"
$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"
$ws.Range("I15").Value = "%"
$ws.Range("J15").Value = "Uninterpretable"
$ws.Range("I20").Value = "aa"
$ws.Range("J20").Value = "Agree/Accept"
$ws.Range("I44").Value = "ba"
$ws.Range("J44").Value = "Appreciation"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "sv"
$ws.Range("J60").Value = "Statement-opinion"
$ws.Range("I63").Value = "sd"
$ws.Range("J63").Value = "Statement-non-opinion"
$ws.Range("I64").Value = "sd"
$ws.Range("J64").Value = "Statement-non-opinion"
$ws.Range("I67").Value = "sv"
$ws.Range("J67").Value = "Statement-opinion"
$ws.Range("I76").Value = "sd"
$ws.Range("J76").Value = "Statement-non-opinion"
$ws.Range("I84").Value = "sv"
$ws.Range("J84").Value = "Statement-opinion"
$ws.Range("I103").Value = "sd"
$ws.Range("J103").Value = "Statement-non-opinion"
$ws.Range("I135").Value = "sv"
$ws.Range("J135").Value = "Statement-opinion"
$ws.Range("I145").Value = "sv"
$ws.Range("J145").Value = "Statement-opinion"
$ws.Range("I155").Value = "sv"
$ws.Range("J155").Value = "Statement-opinion"
$ws.Range("I171").Value = "b"
$ws.Range("J171").Value = "Acknowledge (Backchannel)"
$ws.Range("I172").Value = "b"
$ws.Range("J172").Value = "Acknowledge (Backchannel)"
$ws.Range("I182").Value = "ba"
$ws.Range("J182").Value = "Appreciation"
$ws.Range("I192").Value = "sd"
$ws.Range("J192").Value = "Statement-non-opinion"
$ws.Range("I224").Value = "%"
$ws.Range("J224").Value = "Uninterpretable"
$ws.Range("I225").Value = "%"
$ws.Range("J225").Value = "Uninterpretable"
$ws.Range("I229").Value = "aa"
$ws.Range("J229").Value = "Agree/Accept"
$ws.Range("I250").Value = "%"
$ws.Range("J250").Value = "Uninterpretable"
$ws.Range("I254").Value = "aa"
$ws.Range("J254").Value = "Agree/Accept"
$ws.Range("I267").Value = "aa"
$ws.Range("J267").Value = "Agree/Accept"
$ws.Range("I268").Value = "aa"
$ws.Range("J268").Value = "Agree/Accept"
$ws.Range("I271").Value = "aa"
$ws.Range("J271").Value = "Agree/Accept"
$ws.Range("I273").Value = "sv"
$ws.Range("J273").Value = "Statement-opinion"
$ws.Range("I276").Value = "ba"
$ws.Range("J276").Value = "Appreciation"
$ws.Range("I277").Value = "sv"
$ws.Range("J277").Value = "Statement-opinion"
$ws.Range("I288").Value = "aa"
$ws.Range("J288").Value = "Agree/Accept"
$ws.Range("I293").Value = "sd"
$ws.Range("J293").Value = "Statement-non-opinion"
$ws.Range("I300").Value = "b"
$ws.Range("J300").Value = "Acknowledge (Backchannel)"
$ws.Range("I301").Value = "aa"
$ws.Range("J301").Value = "Agree/Accept"
$ws.Range("I304").Value = "sv"
$ws.Range("J304").Value = "Statement-opinion"
$ws.Range("I312").Value = "sv"
$ws.Range("J312").Value = "Statement-opinion"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I329").Value = "sv"
$ws.Range("J329").Value = "Statement-opinion"
$ws.Range("I333").Value = "sv"
$ws.Range("J333").Value = "Statement-opinion"
$ws.Range("I334").Value = "sv"
$ws.Range("J334").Value = "Statement-opinion"
$ws.Range("I341").Value = "aa"
$ws.Range("J341").Value = "Agree/Accept"
$ws.Range("I369").Value = "sv"
$ws.Range("J369").Value = "Statement-opinion"
$ws.Range("I384").Value = "sd"
$ws.Range("J384").Value = "Statement-non-opinion"
$ws.Range("I386").Value = "aa"
$ws.Range("J386").Value = "Agree/Accept"
$ws.Range("I387").Value = "aa"
$ws.Range("J387").Value = "Agree/Accept"
$ws.Range("I388").Value = "aa"
$ws.Range("J388").Value = "Agree/Accept"
$ws.Range("I427").Value = "sv"
$ws.Range("J427").Value = "Statement-opinion"
$ws.Range("I429").Value = "aa"
$ws.Range("J429").Value = "Agree/Accept"
$ws.Range("I430").Value = "aa"
$ws.Range("J430").Value = "Agree/Accept"
$ws.Range("I440").Value = "sv"
$ws.Range("J440").Value = "Statement-opinion"
$ws.Range("I453").Value = "aa"
$ws.Range("J453").Value = "Agree/Accept"
$ws.Range("I454").Value = "aa"
$ws.Range("J454").Value = "Agree/Accept"
$ws.Range("I455").Value = "ba"
$ws.Range("J455").Value = "Appreciation"
$ws.Range("I456").Value = "sv"
$ws.Range("J456").Value = "Statement-opinion"
